$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column BM (17-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("BM1").Value = "17-aug"

$ws1.Range("BM2").Value  = 67.88
$ws1.Range("BM3").Value  = 60.56
$ws1.Range("BM4").Value  = 53.86
$ws1.Range("BM5").Value  = 41.73
$ws1.Range("BM6").Value  = 37.24
$ws1.Range("BM7").Value  = 27.57
$ws1.Range("BM8").Value  = 29.12
$ws1.Range("BM9").Value  = 17.57
$ws1.Range("BM10").Value = 30.65
$ws1.Range("BM11").Value = 23.68
$ws1.Range("BM12").Value = 9.699999999999999
$ws1.Range("BM13").Value = 0.65
$ws1.Range("BM14").Value = 0.65
$ws1.Range("BM15").Value = 0
$ws1.Range("BM16").Value = 0
$ws1.Range("BM17").Value = 0.65
$ws1.Range("BM18").Value = 5.13
$ws1.Range("BM19").Value = 16.72
$ws1.Range("BM20").Value = 51.45
$ws1.Range("BM21").Value = 94.03
$ws1.Range("BM22").Value = 108.6
$ws1.Range("BM23").Value = 108.05
$ws1.Range("BM24").Value = 101.82
$ws1.Range("BM25").Value = 95.84

# Match header style (bold/centered/bordered) used by the rest of row 1
$ws1.Range("BL1").Copy()
$ws1.Range("BM1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Sheet "Gaz": append new day row 62 ---
# The date-like text "2025-08-15" must stay plain text (like every other
# date cell in this column), so format the cell as Text first to stop
# Excel from auto-converting it into a real date serial, then drop back
# to the default "Normal" style so no stray formatting is left behind.
$ws2 = $wb.Worksheets.Item("Gaz")
$a62 = $ws2.Range("A62")
$a62.NumberFormat = "@"
$a62.Value = "2025-08-15"
$a62.Style = "Normal"
$ws2.Range("B62").Value = 29.825

# --- Sheet "CO2": append new day row 62 ---
$ws3 = $wb.Worksheets.Item("CO2")
$a62b = $ws3.Range("A62")
$a62b.NumberFormat = "@"
$a62b.Value = "2025-08-15"
$a62b.Style = "Normal"
$ws3.Range("B62").Value = 69.95
